# Generate Report for Handback
#
# - Overview sheet: the 1174680b-... row's status (zh-cn + de-de columns)
#   flips from "Ready for handoff" to "Handback transform failed".
# - zh-cn / de-de sheets: the same row's "Error Detail" column gets a
#   diagnostic message explaining the handback/handoff file name mismatch,
#   and the "Error Detail" column is widened to fit the longer text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview: update status for the 1174680b-... row (both locale columns) ---
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# --- zh-cn / de-de: the "Status" column (C3) for the same row shares the very
#     same underlying string as the Overview status, so it flips too ---
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# --- zh-cn: Error Detail message for the 1174680b-... row ---
$wsZhCn.Range("P3").Value = "Handback file name: mret2klu.ifa is different with handoff file name: 1174680b-fb5c-4ec0-89fe-ac75eb728d59.ed5730527e8aef17ac3d8dfc9276fe47cf9491c4.zh-cn."

# --- de-de: Error Detail message for the 1174680b-... row ---
$wsDeDe.Range("P3").Value = "Handback file name: mret2klu.ifa is different with handoff file name: 1174680b-fb5c-4ec0-89fe-ac75eb728d59.ed5730527e8aef17ac3d8dfc9276fe47cf9491c4.de-de."

# --- Widen the "Error Detail" column (P, index 16) to fit the new text ---
# Reuse column A's width (already 40) so the raw stored width matches exactly.
$zhWidth = $wsZhCn.Columns.Item(1).ColumnWidth
$deWidth = $wsDeDe.Columns.Item(1).ColumnWidth
$wsZhCn.Columns.Item(16).ColumnWidth = $zhWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $deWidth
